$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, [string]$value)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

# Row 2
$ws.Cells.Item(2, 4).Value = '29.018.08'
$ws.Cells.Item(2, 5).Value = '  -0.53%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '1.840.03'
$ws.Cells.Item(3, 5).Value = '  -1.00%  '

# Row 4
Set-TextValue $ws.Cells.Item(4, 4) '0.9979'
$ws.Cells.Item(4, 5).Value = '  -0.21%  '

# Row 5
Set-TextValue $ws.Cells.Item(5, 4) '245.45'
$ws.Cells.Item(5, 5).Value = '  +1.35%  '

# Row 6
Set-TextValue $ws.Cells.Item(6, 4) '0.6972'
$ws.Cells.Item(6, 5).Value = '  -0.56%  '

# Row 7
Set-TextValue $ws.Cells.Item(7, 4) '0.9984'
$ws.Cells.Item(7, 5).Value = '  -0.20%  '

# Row 8
Set-TextValue $ws.Cells.Item(8, 4) '0.07711'
$ws.Cells.Item(8, 5).Value = '  -1.14%  '

# Row 9
Set-TextValue $ws.Cells.Item(9, 4) '0.3054'
$ws.Cells.Item(9, 5).Value = '  -1.69%  '

# Row 10
Set-TextValue $ws.Cells.Item(10, 4) '23.47'
$ws.Cells.Item(10, 5).Value = '  -1.77%  '

# Row 11
Set-TextValue $ws.Cells.Item(11, 4) '0.07828'
$ws.Cells.Item(11, 5).Value = '  +0.40%  '

# Row 12
Set-TextValue $ws.Cells.Item(12, 4) '92.87'
$ws.Cells.Item(12, 5).Value = '  +0.25%  '

# Row 13
$ws.Cells.Item(13, 4).Value = '1.834.93'
$ws.Cells.Item(13, 5).Value = '  -1.43%  '

# Row 14
Set-TextValue $ws.Cells.Item(14, 4) '5.113'
$ws.Cells.Item(14, 5).Value = '  -0.24%  '

# Row 15
Set-TextValue $ws.Cells.Item(15, 4) '0.6841'
$ws.Cells.Item(15, 5).Value = '  -1.04%  '

# Row 16
Set-TextValue $ws.Cells.Item(16, 4) '6.592'
$ws.Cells.Item(16, 5).Value = '  +0.63%  '

# Row 17
Set-TextValue $ws.Cells.Item(17, 4) '0.000008275'
$ws.Cells.Item(17, 5).Value = '  -2.05%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '28.977.46'
$ws.Cells.Item(18, 5).Value = '  -0.84%  '

# Row 19
Set-TextValue $ws.Cells.Item(19, 4) '242.16'
$ws.Cells.Item(19, 5).Value = '  -3.17%  '

# Row 20
$ws.Cells.Item(20, 4).Value = '2.076.79'
$ws.Cells.Item(20, 5).Value = '  -1.80%  '

# Row 21
Set-TextValue $ws.Cells.Item(21, 4) '12.71'
$ws.Cells.Item(21, 5).Value = '  -1.67%  '

# Row 22
Set-TextValue $ws.Cells.Item(22, 4) '0.9991'
$ws.Cells.Item(22, 5).Value = '  -0.07%  '

# Row 23
Set-TextValue $ws.Cells.Item(23, 4) '7.482'
$ws.Cells.Item(23, 5).Value = '  -1.55%  '

# Row 24
Set-TextValue $ws.Cells.Item(24, 4) '0.9987'
$ws.Cells.Item(24, 5).Value = '  -0.15%  '

# Row 25
Set-TextValue $ws.Cells.Item(25, 4) '0.1507'
$ws.Cells.Item(25, 5).Value = '  -1.65%  '

# Row 26
Set-TextValue $ws.Cells.Item(26, 4) '158.59'
$ws.Cells.Item(26, 5).Value = '  -1.35%  '

# Row 27
Set-TextValue $ws.Cells.Item(27, 4) '8.800'
$ws.Cells.Item(27, 5).Value = '  -1.21%  '

# Row 28
Set-TextValue $ws.Cells.Item(28, 4) '18.21'
$ws.Cells.Item(28, 5).Value = '  -1.95%  '

# Row 29
Set-TextValue $ws.Cells.Item(29, 4) '1.540'
$ws.Cells.Item(29, 5).Value = '  -1.86%  '

# Row 30
Set-TextValue $ws.Cells.Item(30, 4) '4.224'
$ws.Cells.Item(30, 5).Value = '  -1.13%  '

# Row 31
Set-TextValue $ws.Cells.Item(31, 4) '4.174'
$ws.Cells.Item(31, 5).Value = '  -1.92%  '

# Row 32
Set-TextValue $ws.Cells.Item(32, 4) '1.194'

# Row 33
Set-TextValue $ws.Cells.Item(33, 4) '0.05108'
$ws.Cells.Item(33, 5).Value = '  -2.49%  '

# Row 34
Set-TextValue $ws.Cells.Item(34, 4) '0.7810'
$ws.Cells.Item(34, 5).Value = '  +3.14%  '

# Row 35
Set-TextValue $ws.Cells.Item(35, 4) '1.860'
$ws.Cells.Item(35, 5).Value = '  -0.91%  '

# Row 36
Set-TextValue $ws.Cells.Item(36, 4) '1.146'
$ws.Cells.Item(36, 5).Value = '  -2.57%  '

# Row 37
Set-TextValue $ws.Cells.Item(37, 4) '2.695'
$ws.Cells.Item(37, 5).Value = '  -0.48%  '

# Row 38
$ws.Cells.Item(38, 4).Value = '1.295.96'
$ws.Cells.Item(38, 5).Value = '  +6.18%  '

# Row 39
Set-TextValue $ws.Cells.Item(39, 4) '0.01862'
$ws.Cells.Item(39, 5).Value = '  -0.12%  '

# Row 40
Set-TextValue $ws.Cells.Item(40, 4) '2.699'
$ws.Cells.Item(40, 5).Value = '  -0.85%  '

# Row 41
Set-TextValue $ws.Cells.Item(41, 4) '0.9516'
$ws.Cells.Item(41, 5).Value = '  +5.77%  '

# Row 42
Set-TextValue $ws.Cells.Item(42, 4) '6.145'
$ws.Cells.Item(42, 5).Value = '  +5.40%  '

# Row 43
Set-TextValue $ws.Cells.Item(43, 4) '107.78'
$ws.Cells.Item(43, 5).Value = '  -2.58%  '

# Row 44
Set-TextValue $ws.Cells.Item(44, 4) '0.9988'
$ws.Cells.Item(44, 5).Value = '  -0.08%  '

# Row 45
Set-TextValue $ws.Cells.Item(45, 4) '9.685'
$ws.Cells.Item(45, 5).Value = '  +1.54%  '

# Row 46
$ws.Cells.Item(46, 2).Value = 'Mantle'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue $ws.Cells.Item(46, 4) '0.5166'
$ws.Cells.Item(46, 5).Value = '  -0.34%  '

# Row 47
$ws.Cells.Item(47, 2).Value = 'RocketPoolETH'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Cells.Item(47, 4).Value = '1.977.48'
$ws.Cells.Item(47, 5).Value = '  -1.60%  '

# Row 48
$ws.Cells.Item(48, 2).Value = 'Aave'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Cells.Item(48, 4) '64.01'
$ws.Cells.Item(48, 5).Value = '  -5.02%  '

# Row 49
$ws.Cells.Item(49, 2).Value = 'RenderToken'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Cells.Item(49, 4) '1.757'
$ws.Cells.Item(49, 5).Value = '  -0.65%  '

# Row 50
$ws.Cells.Item(50, 2).Value = 'BabyDogeCoin'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue $ws.Cells.Item(50, 4) '0.00000000119'
$ws.Cells.Item(50, 5).Value = '  -1.96%  '

# Row 51
Set-TextValue $ws.Cells.Item(51, 4) '6.984'
$ws.Cells.Item(51, 5).Value = '  -0.65%  '
